$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.819.39"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "2.083.74"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'233.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").Value = "'59.14"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.52%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +2.37%  "
$ws.Range("D10").Value = "'0.0790"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("E11").Value = "  +2.74%  "
$ws.Range("D12").Value = "2.389.07"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Value = "'14.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.72%  "
$ws.Range("D14").Value = "'21.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("D15").Value = "'0.776"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.06%  "
$ws.Range("D16").Value = "'5.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("D17").Value = "2.075.77"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").Value = "37.728.80"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").Value = "'6.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").Value = "'71.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("D21").Value = "0.0₃0848"
$ws.Range("E21").Value = "  +3.51%  "
$ws.Range("D22").Value = "'228.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("D26").Value = "'9.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.95%  "
$ws.Range("D27").Value = "'171.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("E28").Value = "  -1.66%  "
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("D30").Value = "'19.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("E31").Value = "  +2.37%  "
$ws.Range("E32").Value = "  +2.84%  "
$ws.Range("D33").Value = "'0.0636"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.98%  "
$ws.Range("E34").Value = "  +1.65%  "
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").Value = "'3.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.74%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("E40").Value = "  -1.64%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'17.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.03%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'99.16"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.35%  "
$ws.Range("E43").Value = "  +2.88%  "
$ws.Range("D44").Value = "'2.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.12%  "
$ws.Range("D45").Value = "1.452.45"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("E48").Value = "  +1.73%  "
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.274.50"
$ws.Range("E51").Value = "  +0.35%  "
